$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day-journal entry: copy formatting (date style) of the row above,
# then fill in the new row's data.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 44460
$ws.Range("B11").Value = 40
$ws.Range("C11").Value = "Mieleenpalautus :D Testaa boxscoren fetchaus"

# Leave selection where the author left it
$ws.Range("C14").Select()
